$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "3P" row (row 9) -- its column no longer exists in the new dictionary.
$ws.Rows.Item(9).Delete()

# Make room for a new "2PT" row right after "3P/Game" (now row 9).
$ws.Rows.Item(10).Insert()

# Give the newly inserted row the same look as the surrounding data rows.
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)

# Refill the data-dictionary rows 9-13 with the updated column list.
$ws.Range("A9").Value = "3P/Game"
$ws.Range("B9").Value = "3Pts made per game"
$ws.Range("C9").Value = "Number"
$ws.Range("D9").Value = "Determines 3pt value"

$ws.Range("A10").Value = "2PT"
$ws.Range("B10").Value = "2Pts per game"
$ws.Range("C10").Value = "Number"
$ws.Range("D10").Value = "points from 2pt range"

$ws.Range("A11").Value = "FT"
$ws.Range("B11").Value = "3Pt attempts per game"
$ws.Range("C11").Value = "Number"
$ws.Range("D11").Value = "scoring component"

$ws.Range("A12").Value = "TRB"
$ws.Range("B12").Value = "Total rebound per game"
$ws.Range("C12").Value = "Percentage"
$ws.Range("D12").Value = "heavily impacts WS"

$ws.Range("A13").Value = "AST"
$ws.Range("B13").Value = "Assists per game"
$ws.Range("C13").Value = "Number"
$ws.Range("D13").Value = "Can impact WS"

# Stamp an update note next to the title banner.
$ws.Range("F1").Value = "Updated 12/17/2025"

# Restore the cursor to where the author last left it.
$ws.Range("H12").Select()
